$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row's timestamp-column cell got copied down over the two
# data rows below it, overwriting the timestamp values that used to be
# there with the header text itself.
$headerValue = $ws.Range("A1").Value2
$ws.Range("A2").Value = $headerValue
$ws.Range("A3").Value = $headerValue

# Keep the sheet's frozen header row intact (re-establish it, since it
# otherwise gets lost when the workbook is rewritten) and leave the
# selection where it ended up.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D7").Select() | Out-Null
